# Update codes for recalibration
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # The cell's Range includes a trailing end-of-cell mark; trim the last
    # character (position) so only the visible run text is replaced, leaving
    # the cell mark / paragraph mark and its formatting untouched.
    $target = $d.Range($rng.Start, $rng.End - 1)
    $target.Text = $newText
}

# Resize the third column of the table (tblGrid width 3015 -> 2893 twips
# = 150.75pt -> 144.65pt).
$t.Columns.Item(3).Width = 144.65

# 2-year horizon block
Set-CellText 5 3 "1.93%"
Set-CellText 7 3 "1.42 (1.31 to 1.53)"
Set-CellText 8 3 "0.8% (0.61% to 1%)"
Set-CellText 9 3 "-0.07 (-0.21 to 0.07)"
Set-CellText 10 3 "-0.38 (-0.43 to -0.33)"
Set-CellText 12 3 "0.9 (0.89 to 0.91)"

# 5-year horizon block
Set-CellText 18 3 "5.36%"
Set-CellText 20 3 "0.89 (0.84 to 0.94)"
Set-CellText 21 3 "-0.6% (-0.88% to -0.32%)"
Set-CellText 22 3 "-0.55 (-0.66 to -0.44)"
Set-CellText 23 3 "-0.38 (-0.42 to -0.33)"
Set-CellText 25 3 "0.89 (0.88 to 0.9)"

Write-Output "All replacements applied."
